# Add a new "test3" entry to the sheet, just above the existing E19 ("test2") row.
# This introduces a new shared string ("test3") and a new row (row 18, column E)
# holding it, matching the author's intent of adding another data row for the
# upcoming column-model work described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E18").Value = "test3"
